# Insert a new inventory row ("BIVATRACIN 150MG POWDER SPRAY") right before the
# existing "BOBAI SUNSCREEN ..." row (row 10), shift everything else down by
# one row, renumber the running index column (A), refresh the grand total and
# bump the generated-at timestamp shown in the footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at position 10 (pushes old row 10.. down by one) ---
$ws.Rows.Item(10).Insert()

# --- 2. Populate the new row with the new item's data (same shape/format as
#        every other item row: merged A:B, C:G, H:K, L:M, N:O, single P/Q). ---
$ws.Range("A10:B10").Merge()
$ws.Range("C10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("N10:O10").Merge()

$ws.Cells.Item(10, 1).Value = 4
$ws.Cells.Item(10, 3).Value = "BIVATRACIN 150MG POWDER SPRAY"
$ws.Cells.Item(10, 8).Value = "1:0"
$ws.Cells.Item(10, 14).Value = "83.00"
$ws.Cells.Item(10, 17).Value = "1:0"

# Columns L and P carry numeric-looking custom number formats
# (`#,##0.##;...` / `0.00`), so an auto-typed `.Value` assignment of a
# numeric-looking string would silently become a real number instead of the
# text shared-string the workbook actually stores. Force text entry, then
# restore the original display format.
$fmtL = $ws.Cells.Item(10, 12).NumberFormat
$ws.Cells.Item(10, 12).NumberFormat = "@"
$ws.Cells.Item(10, 12).Value = "1"
$ws.Cells.Item(10, 12).NumberFormat = $fmtL

$fmtP = $ws.Cells.Item(10, 16).NumberFormat
$ws.Cells.Item(10, 16).NumberFormat = "@"
$ws.Cells.Item(10, 16).Value = "83.0000"
$ws.Cells.Item(10, 16).NumberFormat = $fmtP

# --- 3. Renumber the running index in column A for every item row below the
#        insertion point (rows 11..42 now hold items 5..36). ---
for ($r = 11; $r -le 42; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 6
}

# --- Row heights are keyed to the row position itself in this sheet (every
#     row has its own fixed custom height regardless of which item ends up
#     displayed there), so restore the exact original per-row heights after
#     the insert and give the newly created rows their own values. ---
$rowHeights = @{
    10 = 24.75; 11 = 25.5;  12 = 25.5;  13 = 24.75; 14 = 25.5;  15 = 24.75;
    16 = 25.5;  17 = 25.5;  18 = 24.75; 19 = 25.5;  20 = 24.75; 21 = 25.5;
    22 = 25.5;  23 = 24.75; 24 = 25.5;  25 = 24.75; 26 = 25.5;  27 = 25.5;
    28 = 24.75; 29 = 25.5;  30 = 24.75; 31 = 25.5;  32 = 25.5;  33 = 24.75;
    34 = 25.5;  35 = 24.75; 36 = 25.5;  37 = 25.5;  38 = 24.75; 39 = 25.5;
    40 = 24.75; 41 = 25.5;  42 = 25.5;  43 = 24.75; 44 = 16.5
}
foreach ($r in $rowHeights.Keys) {
    $ws.Rows.Item($r).RowHeight = $rowHeights[$r]
}

# --- 4. Update the grand total (old total + new item's price) now on row 43. ---
$ws.Cells.Item(43, 16).Value = 2300.29

# --- 5. Bump the "generated at" timestamp in the footer (now row 44). ---
$ws.Cells.Item(44, 1).Value = "Tuesday, 30 September, 2025 4:08 PM"
